# User selects account type  -- tidy up the sheet view/layout before adding
# the API collection, report and presentation artifacts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the view: scroll back to the top-left (drops the stale
# topLeftCell="F4") and move the selection to B4.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B4").Select()

# Widen column A (was auto "best fit" at ~11.57 chars) to a fixed ~16.57 chars.
$ws.Columns.Item(1).ColumnWidth = 15.7

# Give row 3 more breathing room (109.5 -> 147.75 points) to fit the new text.
$ws.Rows.Item(3).RowHeight = 147.75
